$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 13311.111
$ws.Range("J51").Value = 13311.111
$ws.Range("L51").Value = 13311.111
$ws.Range("N51").Value = -14279.111
$ws.Range("H55").Value = 697
$ws.Range("I55").Value = 697
$ws.Range("K55").Value = 697
$ws.Range("M55").Value = -483
$ws.Range("H70").Value = 359999.66
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 359999.66
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 1079998.98
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -1080538.98
$ws.Range("H73").Value = 359999.66
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 359999.66
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 1079998.98
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -1081870.98
$ws.Range("H100").Value = 2048.2
$ws.Range("I100").Value = 2048.2
$ws.Range("K100").Value = 2048.2
$ws.Range("M100").Value = -1507.2
$ws.Range("H111").Value = 13579.223
$ws.Range("I111").Value = 5554.5
$ws.Range("J111").Value = 19999
$ws.Range("K111").Value = 16663.5
$ws.Range("L111").Value = 59997
$ws.Range("M111").Value = -13596.5
$ws.Range("N111").Value = -66131
$ws.Range("H132").Value = 3375.2156
$ws.Range("I132").Value = 2288.5103
$ws.Range("K132").Value = 6865.5309
$ws.Range("M132").Value = -4335.5309
$ws.Range("H138").Value = 5067.8203
$ws.Range("J138").Value = 6349.1113
$ws.Range("L138").Value = 19047.3339
$ws.Range("N138").Value = -29327.3339

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2873.7
$ws.Range("I45").Value = 3042.125
$ws.Range("J45").Value = 2200
$ws.Range("K45").Value = 3042.125
$ws.Range("L45").Value = 2200
$ws.Range("M45").Value = -2665.125
$ws.Range("N45").Value = -2954

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7328.4375
$ws.Range("I105").Value = 7370.154
$ws.Range("J105").Value = 7147.6665
$ws.Range("K105").Value = 7370.154
$ws.Range("L105").Value = 7147.6665
$ws.Range("M105").Value = -5623.154
$ws.Range("N105").Value = -10641.6665

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5886.609
$ws.Range("J31").Value = 4772.1
$ws.Range("L31").Value = 4772.1
$ws.Range("N31").Value = -5362.1
$ws.Range("H34").Value = 5886.609
$ws.Range("J34").Value = 4772.1
$ws.Range("L34").Value = 4772.1
$ws.Range("N34").Value = -5176.1
$ws.Range("H58").Value = 5080.5713
$ws.Range("I58").Value = 5260.6665
$ws.Range("K58").Value = 5260.6665
$ws.Range("M58").Value = -5057.6665
$ws.Range("H134").Value = 2744.0688
$ws.Range("I134").Value = 1093
$ws.Range("K134").Value = 3279
$ws.Range("M134").Value = -744
$ws.Range("H136").Value = 5080.5713
$ws.Range("I136").Value = 5260.6665
$ws.Range("K136").Value = 15781.9995
$ws.Range("M136").Value = -13231.9995
$ws.Range("H141").Value = 75147
$ws.Range("J141").Value = 79882.28999999999
$ws.Range("L141").Value = 79882.28999999999
$ws.Range("N141").Value = -90242.28999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 554.6667
$ws.Range("I50").Value = 554.6667
$ws.Range("K50").Value = 1664.0001
$ws.Range("M50").Value = -1183.0001
$ws.Range("H53").Value = 554.6667
$ws.Range("I53").Value = 554.6667
$ws.Range("K53").Value = 1664.0001
$ws.Range("M53").Value = -1183.0001
$ws.Range("H109").Value = 3625.4666
$ws.Range("I109").Value = 3323.5
$ws.Range("K109").Value = 9970.5
$ws.Range("M109").Value = -8930.5
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H138").Value = 5608.4287
$ws.Range("I138").Value = 3852
$ws.Range("J138").Value = 9999.5
$ws.Range("K138").Value = 11556
$ws.Range("L138").Value = 29998.5
$ws.Range("M138").Value = -6416
$ws.Range("N138").Value = -40278.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25999
$ws.Range("I70").Value = 25999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 25999
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -25729
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 25999
$ws.Range("I73").Value = 25999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 25999
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -25063
$ws.Range("N73").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4933
$ws.Range("J61").Value = 4933
$ws.Range("L61").Value = 4933
$ws.Range("N61").Value = -5337
$ws.Range("H82").Value = 817.1818
$ws.Range("I82").Value = 369.6
$ws.Range("J82").Value = 1190.1666
$ws.Range("K82").Value = 369.6
$ws.Range("L82").Value = 1190.1666
$ws.Range("M82").Value = -8.600000000000023
$ws.Range("N82").Value = -1912.1666
$ws.Range("H85").Value = 817.1818
$ws.Range("I85").Value = 369.6
$ws.Range("J85").Value = 1190.1666
$ws.Range("K85").Value = 369.6
$ws.Range("L85").Value = 1190.1666
$ws.Range("M85").Value = 878.4
$ws.Range("N85").Value = -3686.1666
$ws.Range("H113").Value = 4933
$ws.Range("J113").Value = 4933
$ws.Range("L113").Value = 4933
$ws.Range("N113").Value = -9273
$ws.Range("H132").Value = 2545.7778
$ws.Range("I132").Value = 2058.8572
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 6176.571599999999
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -3646.571599999999
$ws.Range("N132").Value = -17810
$ws.Range("H136").Value = 2554.6667
$ws.Range("I136").Value = 2499
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 7497
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4947
$ws.Range("N136").Value = -14100

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2221.7354
$ws.Range("I122").Value = 1955.25
$ws.Range("K122").Value = 5865.75
$ws.Range("M122").Value = -3415.75
$ws.Range("H136").Value = 2264.2432
$ws.Range("I136").Value = 1299.6562
$ws.Range("K136").Value = 3898.9686
$ws.Range("M136").Value = -1348.9686
